$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them (e.g. drops
# trailing zeros / turns them into floating point numbers) same as it
# would for a human typing into the grid.
$ws.Range('D2').Value = '26.256.36'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '1.646.27'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.18'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0639'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.02'
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D13').Value = '1.873.30'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').Value = '1.612.59'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('E15').Value = '  -2.47%  '
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.56'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').Value = '26.231.80'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '196.01'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.07'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.63'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('E27').Value = '  +1.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.97'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0507'
$ws.Range('E31').Value = '  +2.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('E35').Value = '  +0.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.914'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('D37').Value = '1.136.33'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.555'
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('E42').Value = '  +2.09%  '
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('D45').Value = '1.782.24'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('E47').Value = '  +3.58%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.418'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.69'
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0972'
$ws.Range('E51').Value = '  +0.99%  '
